$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.972.79"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.073.12"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'236.13"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "'608.17"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").Value = "'0.380"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +9.14%  "
$ws.Range("D11").Value = "3.070.00"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "93.674.34"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'0.0000241"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "'33.80"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "3.639.82"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "3.076.32"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'3.56"
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").Value = "'14.30"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").Value = "'5.75"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "'441.50"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "'8.85"
$ws.Range("E23").Value = "  -5.80%  "
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "'8.35"
$ws.Range("E25").Value = "  +6.04%  "
$ws.Range("D26").Value = "'5.52"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").Value = "'84.45"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").Value = "'11.90"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "3.230.94"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("E31").Value = "  +6.95%  "
$ws.Range("E32").Value = "  +5.25%  "
$ws.Range("E33").Value = "  -7.22%  "
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'8.89"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "'7.44"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").Value = "'0.153"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").Value = "'25.39"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'487.24"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").Value = "'3.85"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'0.436"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'3.08"
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("D47").Value = "'161.69"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'0.675"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "'1.82"
$ws.Range("E49").Value = "  -3.99%  "
$ws.Range("D50").Value = "'43.56"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  +0.09%  "
